$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 256; existing rows 256..337 shift down to 257..338
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256 with the new record's data
$ws.Cells.Item(256, 1).Value = 5
$ws.Cells.Item(256, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(256, 3).Value = "Maule"
$ws.Cells.Item(256, 4).Value = 44876
$ws.Cells.Item(256, 5).Value = 7
$ws.Cells.Item(256, 6).Value = 100112009
$ws.Cells.Item(256, 7).Value = "Acelga"
$ws.Cells.Item(256, 8).Value = "Sin especificar"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 500
$ws.Cells.Item(256, 11).Value = 2000
$ws.Cells.Item(256, 12).Value = 2000
$ws.Cells.Item(256, 13).Value = 2000
$ws.Cells.Item(256, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(256, 15).Value = "Región del Maule"
$ws.Cells.Item(256, 16).Value = 500
$ws.Cells.Item(256, 17).Value = 4
$ws.Cells.Item(256, 18).Value = "Hortaliza"

# Match the date-cell format used by the rest of column D (numFmt applied via style index 2)
$ws.Cells.Item(256, 4).NumberFormat = $ws.Cells.Item(257, 4).NumberFormat
